$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Abu4"
$ws.Range("A3").Value = "Pallu4"
$ws.Range("A4").Value = "Zarina4"
$ws.Range("A5").Value = "Tahira4"
$ws.Range("B2").Value = "abu4@gmail.com"
$ws.Range("B3").Value = "pallu4@gmail.com"
$ws.Range("B4").Value = "zarina4@gmail.com"
$ws.Range("B5").Value = "tahira4@gmail.com"

$ws.Range("B5").Select()
